# Append: 2025-11-03 06:28 JST
# Refresh the "ランサーズ" sheet: replace rows 2-4 with newly scraped
# listings, drop the previously-scraped rows 5-10, and narrow a couple of
# columns to fit the new (shorter) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Clear out all existing hyperlinks (F2:F10) up front so we can
#        re-create only the three that survive (F2:F4) further down. ---
$ws.Range("A1").Hyperlinks.Delete()

# --- 2. Drop the old rows 5-10 entirely (shrinks the used range / dimension
#        down to A1:H4). ---
$ws.Rows("5:10").Delete()

# --- 3. Overwrite rows 2-4 with the refreshed listings. ---

# Row 2
$ws.Range("A2").Value = "2025-11-03 06:28:28"
$ws.Range("B2").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G2").Value = 178
$ws.Range("H2").Value = "★bot ◆ツール"

# Row 3
$ws.Range("A3").Value = "2025-11-03 06:28:28"
$ws.Range("B3").Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Range("G3").Value = 135
$ws.Range("H3").Value = "◆ツール,スクレイピング ◇サイト"

# Row 4
$ws.Range("A4").Value = "2025-11-03 06:28:28"
$ws.Range("B4").Value = "【C#開発】競馬JRA-VAN DataLabを介して、過去レース情報を取得し、DBに保存する業務"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5425801"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆開発"

# --- 4. Re-create hyperlinks for the URL column on rows 2-4, restoring the
#        "Hyperlink" cell style that Add() normally sets, but pinned back to
#        the workbook's existing named style so no stray style gets left
#        behind. ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405023")
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5251319")
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5425801")
$ws.Range("F4").Style = "Hyperlink"

# --- 5. Narrow columns B, D and H to match the shorter refreshed content.
#        COM's ColumnWidth is in "characters" and is offset from the raw
#        OOXML <col width> by the fixed 5/6 char padding Excel adds, so we
#        subtract it here to land on the exact stored widths (51/28/19). ---
$padding = 5 / 6
$ws.Columns("B").ColumnWidth = 51 - $padding
$ws.Columns("D").ColumnWidth = 28 - $padding
$ws.Columns("H").ColumnWidth = 19 - $padding
